# Auto-generated edit script: updates "want to go" (F) / price (G) counters
# per the source-repo refresh described in the commit message.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 247
$ws.Range("F3").Value = 232
$ws.Range("F6").Value = 63
$ws.Range("F8").Value = 2200
$ws.Range("F9").Value = 328
$ws.Range("F10").Value = 33
$ws.Range("F11").Value = 429
$ws.Range("F13").Value = 2536
$ws.Range("F15").Value = 1319
$ws.Range("F16").Value = 4653
$ws.Range("G17").Value = 158
$ws.Range("F18").Value = 4987
$ws.Range("F19").Value = 1557
$ws.Range("F20").Value = 2841
$ws.Range("F21").Value = 3236
$ws.Range("F22").Value = 156
$ws.Range("F23").Value = 1530
$ws.Range("F24").Value = 248
$ws.Range("F25").Value = 834
$ws.Range("F26").Value = 99
$ws.Range("F27").Value = 280
$ws.Range("F28").Value = 955
$ws.Range("F29").Value = 1765
$ws.Range("F30").Value = 114
$ws.Range("F31").Value = 271
$ws.Range("F32").Value = 669
$ws.Range("F34").Value = 321
$ws.Range("F35").Value = 398

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 54
$ws.Range("F8").Value = 88
$ws.Range("F10").Value = 20
$ws.Range("F12").Value = 19
$ws.Range("F14").Value = 43
$ws.Range("F16").Value = 51

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 54
$ws.Range("F7").Value = 247
$ws.Range("F8").Value = 232
$ws.Range("F12").Value = 63
$ws.Range("F13").Value = 2200
$ws.Range("F14").Value = 328
$ws.Range("F15").Value = 88
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 429
$ws.Range("F20").Value = 20
$ws.Range("F21").Value = 2536
$ws.Range("F22").Value = 1319
$ws.Range("F24").Value = 19
$ws.Range("F26").Value = 4653
$ws.Range("G27").Value = 158
$ws.Range("F28").Value = 4987
$ws.Range("F29").Value = 1557
$ws.Range("F30").Value = 2841
$ws.Range("F31").Value = 3236
$ws.Range("F32").Value = 156
$ws.Range("F33").Value = 43
$ws.Range("F35").Value = 1530
$ws.Range("F36").Value = 51
$ws.Range("F37").Value = 248
$ws.Range("F38").Value = 834
$ws.Range("F39").Value = 99
$ws.Range("F40").Value = 280
$ws.Range("F41").Value = 955
$ws.Range("F43").Value = 1765
$ws.Range("F44").Value = 114
$ws.Range("F45").Value = 271
$ws.Range("F46").Value = 669
$ws.Range("F48").Value = 321
$ws.Range("F49").Value = 398
